$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H so the existing H2 formula (and its column)
# shifts to I, leaving room for the new "Casa de discuri" column at H.
$ws.Range("H1").EntireColumn.Insert()

# New column header
$ws.Range("H1").Value = "Casa de discuri"

# Artist name for row 2 changes
$ws.Range("B2").Value = "Baracuda de pe lac"

# The label/distributor that used to live in B2 now lives in H2
$ws.Range("H2").Value = "BY NORSE"

# The lookup formula (previously in H2, referencing B2) now lives in I2 and
# references H2 instead.
$ws.Range("I2").Formula = '=IF(H2="Cooking_vinyl",G2*0.7*1.1,IF(H2="essential",G2*0.7*1.1,IF(H2="one Little Indian",G2*0.7*1.1,IF(H2="Season Of mist",G2*0.7*1.1,IF(H2="frontiers",G2*0.7*1.1,G2*0.75*1.1)))))'
